$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Flip the Runmode column (C) from "Y" to "N" for rows 2-14
# (TestCase_A1 .. TestCase_A13), leaving TestCase_A14 / TestCase_A15
# (rows 15-16) set to run.
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value = "N"
}

# Update the active selection left behind by the author.
$ws.Range("A12").Select()
